$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells that are being updated, to preserve
# them as literal text (matching original inlineStr/shared-string text cells)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values row by row
$ws.Range("D2").Value = '26.628.09'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '1.630.85'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("D5").Value = '213.29'
$ws.Range("E5").Value = '  +0.31%  '
$ws.Range("E6").Value = '  +2.85%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("E8").Value = '  +0.76%  '
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("D10").Value = '19.20'
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").Value = '1.859.98'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.10'
$ws.Range("E13").Value = '  +1.32%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.590.50'
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '63.58'
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.624.69'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").Value = '216.18'
$ws.Range("E19").Value = '  +6.62%  '
$ws.Range("D20").Value = '1.00'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").Value = '6.15'
$ws.Range("E22").Value = '  +1.59%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  +5.47%  '
$ws.Range("D25").Value = '147.55'
$ws.Range("E25").Value = '  +2.03%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("E27").Value = '  +1.03%  '
$ws.Range("E28").Value = '  +3.85%  '
$ws.Range("D29").Value = '15.52'
$ws.Range("E29").Value = '  +2.05%  '
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("E32").Value = '  +3.07%  '
$ws.Range("E33").Value = '  +1.38%  '
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("D36").Value = '1.225.04'
$ws.Range("E36").Value = '  +5.74%  '
$ws.Range("E37").Value = '  +5.56%  '
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("D42").Value = '0.795'
$ws.Range("E42").Value = '  +1.60%  '
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("D44").Value = '1.770.06'
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("D45").Value = '92.82'
$ws.Range("E45").Value = '  +0.83%  '
$ws.Range("E46").Value = '  +2.52%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0104'
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '55.12'
$ws.Range("E48").Value = '  +1.99%  '
$ws.Range("D49").Value = '0.0512'
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("E50").Value = '  +3.53%  '
$ws.Range("E51").Value = '  +0.05%  '
